$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fermenters / Bright Tanks are now part of the brewery - update notes
$ws.Range("G5").Value = "Now part of brewery"
$ws.Range("H5").Value = "N/A"
$ws.Range("G7").Value = "Now part of brewery"
$ws.Range("H7").Value = "N/A"

# New note on Bar Stools
$ws.Range("I9").Value = "but you sit on them!"

# New column J header + summary formula in K1
$ws.Range("J1").Value = "Potential freed ="
$ws.Range("K1").Formula = "=SUM(J:J)-COUNT(J:J)"

# "Potential freed" numbers down column J
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(5, 10).Value = 5
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 3
$ws.Cells.Item(8, 10).Value = 2
$ws.Cells.Item(21, 10).Value = 20
$ws.Cells.Item(22, 10).Value = 4

# Column widths for the new/changed columns
$ws.Columns.Item(7).ColumnWidth = 26.28515625
$ws.Columns.Item(9).ColumnWidth = 19.140625
$ws.Columns.Item(10).ColumnWidth = 14.5703125

$ws.Range("J21").Select()
